$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (plain, no header style applied)
$ws.Range("G1").Value = "Doctor Diagnosis"
$ws.Range("H1").Value = "Treatment Plan"

# New data cells for rows 2 and 3
$ws.Range("G2").Value = "require further followup, require medicine"
$ws.Range("H2").Value = "followup, done"
$ws.Range("G3").Value = "okay, not okay"
$ws.Range("H3").Value = "nil, nil"

# Column G width
$ws.Columns.Item(7).ColumnWidth = 24.33

# Select G4 as the active cell
$ws.Range("G4").Select()

# Zoom level of the sheet view
$excel.ActiveWindow.Zoom = 144
